$wb = $excel.ActiveWorkbook

# 1. Update "sets" sheet: D4 (home_points) 11 -> 12
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Range("D4").Value = 12

# 2. Append new row 75 to "rallies" sheet
$wsRallies = $wb.Worksheets.Item("rallies")

$wsRallies.Range("A75").Value = 74
$wsRallies.Range("B75").Value = 1
$wsRallies.Range("C75").Value = 3
$wsRallies.Range("D75").Value = 12
$wsRallies.Range("E75").Value = "NOS"
$wsRallies.Range("G2").Copy($wsRallies.Range("F75"))
$wsRallies.Range("G75").Value = 3
$wsRallies.Range("H75").Value = "MEIO"
$wsRallies.Range("I75").Value = "PONTO"
$wsRallies.Range("J75").Value = "NOS"
$wsRallies.Range("K75").Value = 12
$wsRallies.Range("L75").Value = 0
$wsRallies.Range("M75").Value = "1 3 m"
$wsRallies.Range("N75").Value = "FRENTE"
$wsRallies.Range("O75").Value = "FRENTE"
$wsRallies.Range("P75").Value = "FRENTE"
